# basicNotation.pptx
#
# Slide 1 has a rectangle ("Rectangle 155") whose first paragraph reads
#   "visibility name = default-value"
# It needs to become
#   "visibility name: type = default-value"
# split across three runs with identical formatting (PowerPoint splits a
# run wherever an in-place text edit is made):
#   "visibility " + "name: type " + "= default-value"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape whose first paragraph starts with "visibility name ="
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text.StartsWith("visibility name =")) {
            $target = $candidate
            break
        }
    }
}

$tr = $target.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Replace the 5 characters "name " (characters 12-16) with "name: type "
$sub = $para1.Characters(12, 5)
$sub.Text = "name: type "
